$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a slightly-rounded Mid Y value on the existing C6 row (row 7) ---
$ws.Range("C7").Value = -46.549999999999997

# --- Insert 9 new blank rows (with formatting copied down) before the old
#     "R3" row (old row 23) so the table can hold the 9 newly-added
#     components while keeping the existing rows' formatting. ---
$ws.Rows.Item(24).Resize(9).Insert()

# --- Now (re)write designator / Mid X / Mid Y / Rotation for every row from
#     14 down to 32 so the whole table ends up alphabetically sorted by
#     designator, with the 9 new parts (C13, C14, D4, D5, R5, R6, R7, U2, Y1)
#     folded in. Layer (column E) stays "top" for all of them, unchanged. ---
$rows = @(
    @(14, "C13", 116.2, -50.549999999999997, -90),
    @(15, "C14", 122.8, -50.549999999999997, -90),
    @(16, "D1", 82.700000000000003, -43.950000000000003, 0),
    @(17, "D2", 82.700000000000003, -46.350000000000001, 180),
    @(18, "D3", 82.700000000000003, -41.549999999999997, 180),
    @(19, "D4", 116.8, -36.994999999999997, 180),
    @(20, "D5", 120.34999999999999, -37.049999999999997, 0),
    @(21, "J1", 97.599999999999994, -39.004398999999999, 180),
    @(22, "JP1", 107.19, -50.740000000000002, 0),
    @(23, "L1", 77.549999999999997, -45.049999999999997, 90),
    @(24, "Q1", 91.200000000000003, -47.200000000000003, -90),
    @(25, "R1", 88, -47.549999999999997, 180),
    @(26, "R2", 100.7, -50.75, 90),
    @(27, "R3", 113.25, -50.600000000000001, 0),
    @(28, "R5", 115.95, -45.350000000000001, 90),
    @(29, "R6", 123.65000000000001, -41.475000000000001, 90),
    @(30, "R7", 123.3, -37.424999999999997, 90),
    @(31, "U2", 119.95, -43.649999999999999, 90),
    @(32, "Y1", 119.5, -50.950000000000003, 180)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = "top"
}

Write-Host "done"
